# Adicionando workflow para publicar imagem no dockerhub
# -> registra a coluna "tratamento" (H) usada na metadado.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cabecalho da nova coluna H
$ws.Range("H1").Value = "tratamento"

# Preenche o tratamento aplicado a cada coluna de origem.
# A ordem das atribuicoes abaixo respeita a ordem em que cada valor de
# texto novo precisa aparecer na tabela de shared strings do workbook.
$ws.Range("H4").Value = "lowercase"
$ws.Range("H2").Value = "n/a"
$ws.Range("H3").Value = "n/a"
$ws.Range("H5").Value = "lowercase"
$ws.Range("H6").Value = "n/a"
$ws.Range("H7").Value = "n/a"
$ws.Range("H8").Value = "n/a"
$ws.Range("H9").Value = "remover caracteres especiais"
$ws.Range("H10").Value = "formatar para Y-m-d"

# Atualiza a celula selecionada ao salvar a planilha
$ws.Range("H3").Select() | Out-Null

# Configuracao de impressao da planilha (tamanho A4, retrato)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
